$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.283.24'
$ws.Range('E2').Value = '  +1.77%  '

$ws.Range('D3').Value = '2.642.84'
$ws.Range('E3').Value = '  +0.61%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = "'605.12"
$ws.Range('E5').Value = '  -0.18%  '

$ws.Range('D6').Value = "'151.26"
$ws.Range('E6').Value = '  +3.12%  '

$ws.Range('D8').Value = "'0.592"
$ws.Range('E8').Value = '  +1.15%  '

$ws.Range('D9').Value = "'0.110"
$ws.Range('E9').Value = '  +1.40%  '

$ws.Range('D10').Value = "'0.392"
$ws.Range('E10').Value = '  +8.20%  '

$ws.Range('D11').Value = "'5.69"
$ws.Range('E11').Value = '  +1.17%  '

$ws.Range('E12').Value = '  -0.67%  '

$ws.Range('D13').Value = "'27.78"
$ws.Range('E13').Value = '  +1.95%  '

$ws.Range('D14').Value = '3.119.26'
$ws.Range('E14').Value = '  +1.01%  '

$ws.Range('D15').Value = '64.139.02'
$ws.Range('E15').Value = '  +1.80%  '

$ws.Range('D16').Value = "'0.0000149"
$ws.Range('E16').Value = '  +3.29%  '

$ws.Range('D17').Value = '2.644.96'
$ws.Range('E17').Value = '  +1.76%  '

$ws.Range('E18').Value = '  +8.47%  '

$ws.Range('E19').Value = '  +4.64%  '

$ws.Range('D20').Value = "'353.22"
$ws.Range('E20').Value = '  +3.86%  '

$ws.Range('D21').Value = "'6.97"
$ws.Range('E21').Value = '  +1.58%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').Value = "'5.74"
$ws.Range('E23').Value = '  +3.20%  '

$ws.Range('D24').Value = "'66.95"
$ws.Range('E24').Value = '  +0.51%  '

$ws.Range('D25').Value = "'1.76"
$ws.Range('E25').Value = '  +14.44%  '

$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = "'1.73"
$ws.Range('E26').Value = '  +6.85%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = "'9.46"
$ws.Range('E27').Value = '  +8.38%  '

$ws.Range('E28').Value = '  +4.98%  '

$ws.Range('D29').Value = "'0.166"
$ws.Range('E29').Value = '  +1.96%  '

$ws.Range('D30').Value = "'543.01"
$ws.Range('E30').Value = '  +1.72%  '

$ws.Range('E31').Value = '  +0.20%  '

$ws.Range('E32').Value = '  +2.17%  '

$ws.Range('D33').Value = '0.0₃0864'
$ws.Range('E33').Value = '  +7.53%  '

$ws.Range('E34').Value = '  +0.94%  '

$ws.Range('D35').Value = "'5.32"
$ws.Range('E35').Value = '  +2.01%  '

$ws.Range('D36').Value = "'168.14"
$ws.Range('E36').Value = '  -0.92%  '

$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = "'2.03"
$ws.Range('E37').Value = '  +8.00%  '

$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = "'0.411"
$ws.Range('E38').Value = '  +2.12%  '

$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('D40').Value = "'19.66"
$ws.Range('E40').Value = '  +3.40%  '

$ws.Range('E41').Value = '  +0.08%  '

$ws.Range('D42').Value = "'167.49"
$ws.Range('E42').Value = '  -1.55%  '

$ws.Range('D43').Value = "'40.40"
$ws.Range('E43').Value = '  +1.83%  '

$ws.Range('E44').Value = '  +5.45%  '

$ws.Range('D45').Value = "'0.0584"
$ws.Range('E45').Value = '  +2.82%  '

$ws.Range('D46').Value = "'21.75"
$ws.Range('E46').Value = '  -2.85%  '

$ws.Range('D47').Value = "'0.631"
$ws.Range('E47').Value = '  +1.15%  '

$ws.Range('E48').Value = '  +14.39%  '

$ws.Range('D49').Value = "'0.0247"
$ws.Range('E49').Value = '  +2.96%  '

$ws.Range('D50').Value = "'0.0968"
$ws.Range('E50').Value = '  +0.78%  '

$ws.Range('D51').Value = "'19.42"
$ws.Range('E51').Value = '  +5.12%  '
